$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 3 earlier employee rows (rows 2-4), keeping header (row 1)
# and the last employee's row, which becomes the new row 2.
$ws.Rows("2:4").Delete()

# The header row's explicit height was sized for wrapped text across the
# original wider column set; re-fit it now that columns have shrunk.
$ws.Rows("1").AutoFit()

# Remove unused columns: Assets (M), Designation (D), CompanyId (A).
# Deleted right-to-left so earlier deletes don't shift later column letters.
$ws.Columns("M").Delete()
$ws.Columns("D").Delete()
$ws.Columns("A").Delete()

# Update the surviving row's Role and Email values.
$ws.Range("E2").Value = "Jthomson@crossleaf.com"
$ws.Range("A2").Value = "Developer"

# Preserve the existing hyperlink-style formatting on E2 before we touch
# the Hyperlinks collection (Hyperlinks.Add() reformats the target cell).
$ws.Range("E2").Copy()
$ws.Range("Z100").PasteSpecial(-4122)

# Drop the stale hyperlinks (formerly pointing at G2:G5) and add the one
# hyperlink that belongs to the remaining row.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:Jthomson@crossleaf.com")

# Restore E2's original formatting, then clean up the scratch cell.
$ws.Range("Z100").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("Z100").Clear()

# Match the saved selection/view state.
$ws.Range("A2").Select() | Out-Null
